{"js": "// Insert two \"Author\" paragraphs (\"Laura Kennedy\" and \"Ben Jarman\")\n// immediately after the Title paragraph (\"Research Data Management: An\n// introduction\"), i.e. before the Date paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,style,text\");\nawait context.sync();\n\n// Find the Title paragraph (first paragraph with the \"Title\" style).\nlet titlePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].style === \"Title\") {\n    titlePara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!titlePara) {\n  throw new Error(\"Could not find Title paragraph\");\n}\n\n// insertParagraph inserts relative to the reference paragraph; inserting\n// \"After\" each time right after the title keeps the authors directly\n// below it and in the requested order (Laura Kennedy, then Ben Jarman).\nconst benPara = titlePara.insertParagraph(\"Ben Jarman\", Word.InsertLocation.after);\nbenPara.style = \"Author\";\n\nconst lauraPara = titlePara.insertParagraph(\"Laura Kennedy\", Word.InsertLocation.after);\nlauraPara.style = \"Author\";\n\nawait context.sync();\n", "ps1": "# Insert two \"Author\" paragraphs (\"Laura Kennedy\" and \"Ben Jarman\")\n# immediately after the Title paragraph (\"Research Data Management: An\n# introduction\"), i.e. before the Date paragraph.\n\n$d = $word.ActiveDocument\n\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$lauraPara = $d.Paragraphs(2)\n$lauraPara.Range.Text = \"Laura Kennedy\"\n$lauraPara.Style = \"Author\"\n\n$lauraPara.Range.InsertParagraphAfter()\n\n$benPara = $d.Paragraphs(3)\n$benPara.Range.Text = \"Ben Jarman\"\n$benPara.Style = \"Author\"\n"}
